$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsRHPF  = $wb.Worksheets.Item("RHPF")

# --- About sheet: update the note text (A12) ---
$wsAbout.Range("A12").Value = "electrolysis that is guaranteed to be supplied by new clean electricity sources."

# --- RHPF sheet: rename the "thermochemical water splitting" pathway to
#     "hydrocarbon partial oxidation" (column F header + row label) ---
$wsRHPF.Range("F1").Value = "hydrocarbon partial oxidation"
$wsRHPF.Range("A6").Value = "hydrocarbon partial oxidation"

# --- RHPF sheet: update the fraction values ---
# Row 2 (electrolysis) : 0.05 -> 0 across B:H
$wsRHPF.Range("B2:F2").Value = 0
$wsRHPF.Range("G2:H2").ClearFormats()
$wsRHPF.Range("G2:H2").Value = 0

# Row 3 (natural gas reforming) : 0.95 -> 0 across B:H
$wsRHPF.Range("B3:F3").Value = 0
$wsRHPF.Range("G3:H3").ClearFormats()
$wsRHPF.Range("G3:H3").Value = 0

# Row 7 (electrolysis with guaranteed clean electricity) : 0 -> 1 across B:H
$wsRHPF.Range("B7:F7").Value = 1
$wsRHPF.Range("G7:H7").ClearFormats()
$wsRHPF.Range("G7:H7").Value = 1

# --- Restore selections recorded in the saved file ---
$wsRHPF.Range("F2").Select()
$wsAbout.Range("A13").Select()
$wsAbout.Activate()
